$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.40"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("D3").Value = "'26.93"
$ws.Range("E3").Value = "'-0.11%"
$ws.Range("D4").Value = "'4.695"
$ws.Range("E4").Value = "'1.65%"
$ws.Range("D5").Value = "'0.06049"
$ws.Range("E5").Value = "'2.88%"
$ws.Range("D6").Value = "'6.685"
$ws.Range("E6").Value = "'0.93%"
$ws.Range("D7").Value = "'0.8600"
$ws.Range("E7").Value = "'0.25%"
$ws.Range("D8").Value = "'0.9240"
$ws.Range("E8").Value = "'-1.94%"
$ws.Range("E9").Value = "'-0.96%"
$ws.Range("D10").Value = "'0.05133"
$ws.Range("E10").Value = "'25.51%"
$ws.Range("D11").Value = "'0.07083"
$ws.Range("E11").Value = "'-0.17%"
$ws.Range("D12").Value = "'0.03069"
$ws.Range("E12").Value = "'-3.63%"
$ws.Range("D13").Value = "'0.09139"
$ws.Range("E13").Value = "'-0.23%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("D15").Value = "'0.0006060"
$ws.Range("E15").Value = "'-0.03%"
$ws.Range("D16").Value = "'0.006044"
$ws.Range("E16").Value = "'-3.00%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'-1.26%"
$ws.Range("E18").Value = "'-1.15%"
$ws.Range("E19").Value = "'-0.94%"
$ws.Range("D20").Value = "'0.3128"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("D22").Value = "'4.116"
$ws.Range("E22").Value = "'7.41%"
$ws.Range("D23").Value = "'0.04243"
$ws.Range("E23").Value = "'0.37%"
$ws.Range("E24").Value = "'-0.77%"
$ws.Range("D25").Value = "'0.004022"
$ws.Range("E25").Value = "'-6.24%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.06%"
$ws.Range("E27").Value = "'-21.38%"
$ws.Range("D40").Value = "'0.03857"
$ws.Range("E40").Value = "'0.78%"
$ws.Range("D41").Value = "'0.1116"
$ws.Range("E41").Value = "'1.41%"
$ws.Range("D42").Value = "'0.004044"
$ws.Range("E42").Value = "'-34.88%"
$ws.Range("D43").Value = "'0.01477"
$ws.Range("E43").Value = "'29.38%"
$ws.Range("E44").Value = "'-0.04%"
$ws.Range("D45").Value = "'0.00005189"
$ws.Range("E45").Value = "'-5.05%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.1353"
$ws.Range("E47").Value = "'-42.34%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.05453"
$ws.Range("E48").Value = "'6.45%"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'-0.04%"
